$d = $word.ActiveDocument

# Locate the run of text that needs to be split so that the name/email
# portion can be highlighted (redacted) in black.
$target = "Professor Michael Frank at mcfrank@stanford.edu."

$findRange = $d.Content
$found = $findRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $nameRange = $d.Range($findRange.Start, $findRange.End)
    $nameRange.Font.HighlightColorIndex = 1  # wdBlack
} else {
    Write-Host "Target text not found!"
}
